# Generate Report for Handback
#
# Updates the handback-status report with the freshly generated
# handback timestamps for the second source file
# (9dced597-287e-4541-a4e2-02146d0d2a6d.md):
#   - Overview sheet:      "de-de" column (latest HO xliff generate date)
#   - zh-cn sheet:          Correspond Handback DateTime
#   - de-de sheet:          Correspond Handoff Datetime
#   - de-de sheet:          Correspond Handback DateTime

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G3 - de-de latest HO xliff generate date for 9dced597...md
$wsOverview.Range("G3").Value = "2016-09-01 09:02:26"

# zh-cn!K3 - Correspond Handback DateTime for 9dced597...md
$wsZhCn.Range("K3").Value = "2016-09-01 09:02:43"

# de-de!H3 - Correspond Handoff Datetime for 9dced597...md
$wsDeDe.Range("H3").Value = "2016-09-01 09:02:26"

# de-de!K3 - Correspond Handback DateTime for 9dced597...md
$wsDeDe.Range("K3").Value = "2016-09-01 09:02:50"
